# fix excel export and tests
#
# 1. Reaction sheet: units "umole/min" -> "umole/minute"
# 2. UnitDef sheet: rebuilt as the merged "action" table (setNS + defineUnit
#    rows, a new "type" column) and renamed to "undefined"
# 3. The old standalone "undefined" sheet (setNS/action table) is removed -
#    its single data row now lives as row 2 of the rebuilt sheet above.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Reaction!F2 units fix
# ---------------------------------------------------------------------
$wsReaction = $wb.Worksheets.Item("Reaction")
$wsReaction.Range("F2").Value = "umole/minute"

# ---------------------------------------------------------------------
# 2) Drop the old "undefined" sheet (sheet6) - its one data row gets
#    folded into the rebuilt UnitDef sheet as row 2 below.
# ---------------------------------------------------------------------
$wsOldUndefined = $wb.Worksheets.Item("undefined")
$null = $wsOldUndefined.Delete()

# ---------------------------------------------------------------------
# 3) Rebuild the UnitDef sheet's data, then rename it to "undefined"
# ---------------------------------------------------------------------
$wsUnitDef = $wb.Worksheets.Item("UnitDef")

# wipe everything first so no stale cells (e.g. old C column "mm" values,
# or the now-gone row 50) survive
$wsUnitDef.Range("A1:F50").ClearContents()

# header row
$wsUnitDef.Cells.Item(1,1).Value = "on"
$wsUnitDef.Cells.Item(1,2).Value = "action"
$wsUnitDef.Cells.Item(1,3).Value = "space"
$wsUnitDef.Cells.Item(1,4).Value = "id"
$wsUnitDef.Cells.Item(1,5).Value = "units"
$wsUnitDef.Cells.Item(1,6).Value = "type"

# row 2: the setNS/"concrete" action row (formerly the "undefined" sheet)
$wsUnitDef.Cells.Item(2,1).Value = 1
$wsUnitDef.Cells.Item(2,2).Value = "setNS"
$wsUnitDef.Cells.Item(2,3).Value = "mm"
$wsUnitDef.Cells.Item(2,6).Value = "concrete"

# rows 3-50: defineUnit rows (id in D, units-definition in E)
$units = @(
    @("fmole", "(1e-15 mole)"),
    @("pmole", "(1e-12 mole)"),
    @("nmole", "(1e-9 mole)"),
    @("umole", "(1e-6 mole)"),
    @("mmole", "(1e-3 mole)"),
    @("fM", "(1e-15 mole)/litre"),
    @("pM", "(1e-12 mole)/litre"),
    @("nM", "(1e-9 mole)/litre"),
    @("uM", "(1e-6 mole)/litre"),
    @("mM", "(1e-3 mole)/litre"),
    @("M", "mole/litre"),
    @("kM", "(1e+3 mole)/litre"),
    @("fL", "(1e-15 litre)"),
    @("pL", "(1e-12 litre)"),
    @("nL", "(1e-9 litre)"),
    @("uL", "(1e-6 litre)"),
    @("mL", "(1e-3 litre)"),
    @("dL", "(1e-1 litre)"),
    @("L", "litre"),
    @("fs", "(1e-15 second)"),
    @("ps", "(1e-12 second)"),
    @("ns", "(1e-9 second)"),
    @("us", "(1e-6 second)"),
    @("ms", "(1e-3 second)"),
    @("s", "second"),
    @("h", "hour"),
    @("week", "(7e+0 day)"),
    @("fg", "(1e-18 kilogram)"),
    @("pg", "(1e-15 kilogram)"),
    @("ng", "(1e-12 kilogram)"),
    @("ug", "(1e-9 kilogram)"),
    @("mg", "(1e-6 kilogram)"),
    @("g", "(1e-3 kilogram)"),
    @("kg", "kilogram"),
    @("kat", "katal"),
    @("cell", "item"),
    @("kcell", "(1e+3 item)"),
    @("cal", "(4.1868e+0 joule)"),
    @("kcal", "(4.1868e+3 joule)"),
    @("fm", "(1e-15 metre)"),
    @("pm", "(1e-12 metre)"),
    @("nm", "(1e-9 metre)"),
    @("um", "(1e-6 metre)"),
    @("mm", "(1e-13 metre)"),
    @("cm", "(1e-2 metre)"),
    @("m", "metre"),
    @("UL", "dimensionless"),
    @("percent", "(1e-2 dimensionless)")
)

$row = 3
foreach ($u in $units) {
    $wsUnitDef.Cells.Item($row,1).Value = 1
    $wsUnitDef.Cells.Item($row,2).Value = "defineUnit"
    $wsUnitDef.Cells.Item($row,4).Value = $u[0]
    $wsUnitDef.Cells.Item($row,5).Value = $u[1]
    $row = $row + 1
}

$wsUnitDef.Name = "undefined"
